$wb = $excel.ActiveWorkbook

# Add the new "contact" sheet after the last existing sheet so it lands
# at the end of the tab strip (and becomes the active sheet), matching
# the target workbook layout.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "contact"

$ws.Range("A1:E1").Merge()
$ws.Range("A1").Value = "Example author"
$ws.Range("A1:E1").HorizontalAlignment = -4108  # xlCenter

$ws.Range("A2").Value = "Marian"
$ws.Range("B2").Value = "Marianacki"
$ws.Range("A3").Value = 123
$ws.Range("A4").Value = "Pila-wojenna"

$ws.Range("A5").Select() | Out-Null
